$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Shift the existing rows 15-41 down to 19-45 to make room for 4 new
#    "Linked covariates" rows describing elevation / climate / temperature /
#    precipitation covariates.
# ---------------------------------------------------------------------------
$srcRange = $ws.Range("A15:D41")
$shiftedValues = $srcRange.Value2
$dstRange = $ws.Range("A19:D45")
$dstRange.Value2 = $shiftedValues

# ---------------------------------------------------------------------------
# 2) Fix up formatting so that:
#    - the new bottom-most row (45) gets the "last row" bottom-border style
#      (copied from the row that used to be last, row 41)
#    - all rows from 15 to 44 use the regular "middle" row style (copied
#      from a plain row, e.g. row 13)
#    NOTE: order matters - copy the bottom-border style to row 45 BEFORE
#    overwriting rows 15-44 (which include the old row 41) with the plain
#    style, otherwise the source style would already be gone.
# ---------------------------------------------------------------------------
$lastRowStyleSrc = $ws.Range("A41:D41")
$lastRowStyleSrc.Copy()
$ws.Range("A45:D45").PasteSpecial(-4122)

$plainRowStyleSrc = $ws.Range("A13:D13")
$plainRowStyleSrc.Copy()
$ws.Range("A15:D44").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Write the 4 new "Linked covariates" rows (15-18).
# ---------------------------------------------------------------------------
$newRows = New-Object 'object[,]' 4,4

$newRows[0,0] = "Linked covariates"
$newRows[0,1] = "The average elevation estimated within the spatial domain of the Urban Centre, and expressed in metres above sea level (MASL) (EORC & JAXA, 2017)."
$newRows[0,2] = "EL_AV_ALS"
$newRows[0,3] = "city"

$newRows[1,0] = "Linked covariates"
$newRows[1,1] = "Semi-colon separated list of names of Köppen-Geiger climate classes, intersecting with the spatial domain of the Urban Centre (1986-2010) (Rubel et al., 2017)."
$newRows[1,2] = "E_KG_NM_LST"
$newRows[1,3] = "city"

$newRows[2,0] = "Linked covariates"
$newRows[2,1] = "Average temperature calculated from annual average estimates for time interval centred on the year 2015 (the interval spans from 2012 to 2015) within the spatial domain of the Urban Centre, and expressed in Celsius degrees (°C) (Harris et al., 2014)."
$newRows[2,2] = "E_WR_T_14"
$newRows[2,3] = "city"

$newRows[3,0] = "Linked covariates"
$newRows[3,1] = "Average precipitations calculated from annual average estimates for time interval centred on the year 2015 (the interval spans from 2012 to 2015) within the spatial domain of the Urban Centre; and expressed in millimetres (mm), the amount of rain per square meter in one hour) (Harris et al., 2014)."
$newRows[3,2] = "E_WR_P_14"
$newRows[3,3] = "city"

$ws.Range("A15:D18").Value2 = $newRows

# ---------------------------------------------------------------------------
# 4) Update the "Study region" row - variable name changed from "City" to
#    "study_region".
# ---------------------------------------------------------------------------
$ws.Range("C5").Value2 = "study_region"

# ---------------------------------------------------------------------------
# 5) Misc workbook-level metadata / view changes recorded in the diff.
# ---------------------------------------------------------------------------
$ws.Range("B18").Select()

$wnd = $excel.ActiveWindow
$wnd.ScrollRow = 1
